$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)

    $ws.Range("F2").Value = 679
    $ws.Range("F3").Value = 3967
    $ws.Range("F4").Value = 110
    $ws.Range("F5").Value = 735
    $ws.Range("G5").Value = 50
}
